$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("債務")

# Pre-format the "date" column (J) as text so the ISO-looking date string
# "2011-12-26" is stored as a literal string instead of being auto-converted
# to a date serial number by Excel's input parser.
$ws.Range("J1:J4").NumberFormat = "@"

# Row 1: header row (field names)
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Row 2: year 96 debt (房屋貸款 / 第一銀行恆春分行)
$ws.Range("B2").Value = "房屋貸款"
$ws.Range("C2").Value = "潘孟安"
$ws.Range("D2").Value = "第一銀行恆春分行"
$ws.Range("E2").Value = 2617597
$ws.Range("F2").Value = "98年12月24日"
$ws.Range("G2").Value = "房貸"
$ws.Range("H2").Value = "debt"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2011-12-26"
$ws.Range("K2").Value = "潘孟安"
$ws.Range("L2").Value = 1376
$ws.Range("M2").Value = "tmp6a821"
$ws.Range("N2").Value = 96

# Row 3: year 97 debt (房屋貸款 / 彰化銀行車城分行)
$ws.Range("B3").Value = "房屋貸款"
$ws.Range("C3").Value = "潘孟安"
$ws.Range("D3").Value = "彰化銀行車城分行"
$ws.Range("E3").Value = 5000000
$ws.Range("F3").Value = "98年07月06日"
$ws.Range("G3").Value = "房貸"
$ws.Range("H3").Value = "debt"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").Value = "2011-12-26"
$ws.Range("K3").Value = "潘孟安"
$ws.Range("L3").Value = 1376
$ws.Range("M3").Value = "tmp6a821"
$ws.Range("N3").Value = 97

# Row 4: year 98 debt (信用貸款 / 台灣銀行中屏分行)
$ws.Range("B4").Value = "信用貸款"
$ws.Range("C4").Value = "潘孟安"
$ws.Range("D4").Value = "台灣銀行中屏分行"
$ws.Range("E4").Value = 69179
$ws.Range("F4").Value = "98年07月02日"
$ws.Range("G4").Value = "信貸"
$ws.Range("H4").Value = "debt"
$ws.Range("I4").Value = "normal"
$ws.Range("J4").Value = "2011-12-26"
$ws.Range("K4").Value = "潘孟安"
$ws.Range("L4").Value = 1376
$ws.Range("M4").Value = "tmp6a821"
$ws.Range("N4").Value = 98
